$wb = $excel.ActiveWorkbook

# --- Insert the new "3.1" worksheet before "4.1" --------------------------
$target = $wb.Worksheets.Item("4.1")
$ws = $wb.Worksheets.Add($target)
$ws.Name = "3.1"

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "t"
$ws.Range("B1").Value = "beta"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("A1:B1").Borders.Item(9).LineStyle = 1

# --- Data rows ----------------------------------------------------------
$ws.Range("A2").Value = 0.125
$ws.Range("A3").Formula = "=A2+0.125"
$ws.Range("A4").Formula = "=A3+0.125"
$ws.Range("A5").Formula = "=A4+0.125"
$ws.Range("A6").Formula = "=A5+0.125"
$ws.Range("A7").Formula = "=A6+0.125"
$ws.Range("A8").Formula = "=A7+0.125"
$ws.Range("A9").Formula = "=A8+0.125"

$ws.Range("B2").Formula = "=(2.73/((14*A2)^2))^0.25"
$ws.Range("B3").Formula = "=(2.73/((14*A3)^2))^0.25"
$ws.Range("B4").Formula = "=(2.73/((14*A4)^2))^0.25"
$ws.Range("B5").Formula = "=(2.73/((14*A5)^2))^0.25"
$ws.Range("B6").Formula = "=(2.73/((14*A6)^2))^0.25"
$ws.Range("B7").Formula = "=(2.73/((14*A7)^2))^0.25"
$ws.Range("B8").Formula = "=(2.73/((14*A8)^2))^0.25"
$ws.Range("B9").Formula = "=(2.73/((14*A9)^2))^0.25"

$ws.Range("A2:B9").NumberFormat = "0.000"
$ws.Range("A2:B9").HorizontalAlignment = -4108

# --- Column widths (match default 9.140625 used elsewhere) --------------
$ws.Columns.Item("A:B").ColumnWidth = 9.140625

# --- Page layout ----------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- Sheet view: make "3.1" the active/selected tab, mirroring the source -
$ws.Activate()
$ws.Range("L14:M16").Select()

$wb.Save()
